$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 8 (Type_adresse entry) - causes everything below to shift up by one.
$ws.Rows.Item(8).Delete()
